$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Slit1","Robo1","ECs",3,1,0.3766166666666667,1.12985,0.6239986391709988,0.6239986391709988,2,0.6666666666666666,0.1757713333333334,0.5273140000000001,0.009444264308298454,0.009444264308298454,0.06619841365555557,0.5957857229000001,0.00589320807634947,0.00589320807634947),
    @("ECs","Slit1","Robo1","FAPs",3,1,0.3766166666666667,1.12985,0.6239986391709988,0.6239986391709988,3,1,16.98312366666667,50.949371,0.9125100529581165,0.9125100529581165,6.396127424927778,57.56514682435,0.5694050312757207,0.5694050312757207),
    @("ECs","Slit1","Robo1","sCs",3,1,0.3766166666666667,1.12985,0.6239986391709988,0.6239986391709988,3,1,1.452542333333333,4.357627,0.07804568273358503,0.07804568273358505,0.5470516517722223,4.92346486595,0.04870039981892858,0.04870039981892859),
    @("FAPs","Slit1","Robo1","ECs",2,0.6666666666666666,0.2055996666666667,0.616799,0.340648525593692,0.340648525593692,2,0.6666666666666666,0.1757713333333334,0.5273140000000001,0.009444264308298454,0.009444264308298454,0.03613852754288889,0.3252467478860001,0.003217174711938998,0.003217174711938998),
    @("FAPs","Slit1","Robo1","FAPs",2,0.6666666666666666,0.2055996666666667,0.616799,0.340648525593692,0.340648525593692,3,1,16.98312366666667,50.949371,0.9125100529581165,0.9125100529581165,3.491724564825444,31.425521083429,0.3108452041296042,0.3108452041296042),
    @("FAPs","Slit1","Robo1","sCs",2,0.6666666666666666,0.2055996666666667,0.616799,0.340648525593692,0.340648525593692,3,1,1.452542333333333,4.357627,0.07804568273358503,0.07804568273358505,0.2986422195525555,2.687779975973,0.02658614675214881,0.02658614675214881),
    @("sCs","Slit1","Robo1","ECs",2,0.6666666666666666,0.02133733333333333,0.064012,0.0353528352353091,0.03535283523530909,2,0.6666666666666666,0.1757713333333334,0.5273140000000001,0.009444264308298454,0.009444264308298454,0.003750491529777778,0.033754423768,0.0003338815200099857,0.0003338815200099856),
    @("sCs","Slit1","Robo1","FAPs",2,0.6666666666666666,0.02133733333333333,0.064012,0.0353528352353091,0.03535283523530909,3,1,16.98312366666667,50.949371,0.9125100529581165,0.9125100529581165,0.3623745707168889,3.261371136452,0.03225981755279147,0.03225981755279147),
    @("sCs","Slit1","Robo1","sCs",2,0.6666666666666666,0.02133733333333333,0.064012,0.0353528352353091,0.03535283523530909,3,1,1.452542333333333,4.357627,0.07804568273358503,0.07804568273358505,0.03099337994711111,0.278940419524,0.00275913616250764,0.00275913616250764)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
